$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewsPaper")
$ws.Activate()

# 1. Insert a new column before B ("Paper") to hold "Distribution Location"
$ws.Columns("B:B").Insert()

# 2. Insert two new rows before the old row 6 ("New Garia") to hold the new
#    "Garia / Garia Bata Shop / Sanjeev Paul" distribution entry
$ws.Rows("6:7").Insert()

# 3. Populate the new "Garia Bata Shop" distribution row (row 6) and its
#    blank "Anandabazar" placeholder row (row 7), following the pattern
#    used by the other location blocks on this sheet.
$ws.Range("A6").Value = "Garia"
$ws.Range("B6").Value = "Garia Bata Shop"
$ws.Range("C6").Value = "ToI"
$ws.Range("D6").Value = "Sanjeev Paul"
$ws.Range("E6").Value = 2000
$ws.Range("F6").Value = "SUNDAY"
$ws.Range("G6").Value = 0.1875
$ws.Range("G6").NumberFormat = "h:mm AM/PM"
$ws.Range("I6").Value = 9831570813
$ws.Range("J6").Value = "Call him on 5th December for booking on 10th December"
$ws.Range("K6").Value = 45270
$ws.Range("K6").NumberFormat = "mm-dd-yy"

$ws.Range("C7").Value = "Anandabazar"

# 4. New column headers for the appended "Remarks" / "Distribution Date" columns
$ws.Range("B1").Value = "Distribution Location"
$ws.Range("J1").Value = "Remarks"
$ws.Range("K1").Value = "Distribution Date"

# 5. Remarks noted against the existing distribution rows
$ws.Range("J2").Value = "Did not pick up the call"
$ws.Range("J4").Value = "Call Sanjeev Paul"
$ws.Range("J10").Value = "Getting Switched Off"

# 6. Highlight the Garia/GAUTAM PATRA row (now stale) in red
$ws.Range("A4:J5").Interior.Color = 255

Write-Host "Structure changes done"
